# Weekly fruit/vegetable price log: a new daily record was inserted above
# the current row 237 (shifting rows 237:286 down to 238:287, and growing
# the used range from A1:R286 to A1:R287). The new row carries the same
# market/product metadata as the row that used to sit at 237, but with a
# fresh date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 237; everything from 237 down shifts to 238.
$ws.Rows(237).Insert()

# Populate the newly inserted row 237.
$ws.Cells.Item(237, 1).Value = 3
$ws.Cells.Item(237, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(237, 3).Value = "Coquimbo"
$ws.Cells.Item(237, 4).Value = 44637
$ws.Cells.Item(237, 5).Value = 5
$ws.Cells.Item(237, 6).Value = 100112039
$ws.Cells.Item(237, 7).Value = "Ciboulette"
$ws.Cells.Item(237, 8).Value = "Sin especificar"
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 160
$ws.Cells.Item(237, 11).Value = 1500
$ws.Cells.Item(237, 12).Value = 1500
$ws.Cells.Item(237, 13).Value = 1500
$ws.Cells.Item(237, 14).Value = "`$/docena de atados"
$ws.Cells.Item(237, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(237, 16).Value = 500
$ws.Cells.Item(237, 17).Value = 3
$ws.Cells.Item(237, 18).Value = "Hortaliza"
